# #5: insurance, claim, debt, investment done
#
# The "債務" (debt) worksheet only had a bare/incomplete schema (species,
# debtor, owner, register_date, register_reason repeated as both header and
# data). Finish it out like the other sheets: add property_category,
# category, date, legislator_name, legislator_id, source_file, index, and
# give the header row its proper field-name labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("債務")

function Set-TextValue($cell, [string]$text) {
    # Writing a literal like "2012-04-30" through .Value/.Formula gets
    # auto-coerced into a date serial by Excel's type inference. Routing it
    # through a text formula + paste-as-values keeps it a plain string.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# --- Header row (row 1) : give it real field names ------------------------
$ws.Cells.Item(1, 2).Value = "species"
$ws.Cells.Item(1, 3).Value = "debtor"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "total"
$ws.Cells.Item(1, 6).Value = "register_date"
$ws.Cells.Item(1, 7).Value = "register_reason"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Row 2 (index 113) ------------------------------------------------------
$ws.Cells.Item(2, 2).Value = "貸款"
$ws.Cells.Item(2, 6).Value = "91年4月"
$ws.Cells.Item(2, 7).Value = "房屋貸款"
$ws.Cells.Item(2, 8).Value = "debt"
$ws.Cells.Item(2, 9).Value = "normal"
Set-TextValue $ws.Cells.Item(2, 10) "2012-04-30"
$ws.Cells.Item(2, 11).Value = "蕭美琴"
$ws.Cells.Item(2, 12).Value = 981
$ws.Cells.Item(2, 13).Value = "tmpcd9a1"
$ws.Cells.Item(2, 14).Value = 113

# --- Row 3 (index 114) ------------------------------------------------------
$ws.Cells.Item(3, 2).Value = "貸款"
$ws.Cells.Item(3, 7).Value = "房屋貸款"
$ws.Cells.Item(3, 8).Value = "debt"
$ws.Cells.Item(3, 9).Value = "normal"
Set-TextValue $ws.Cells.Item(3, 10) "2012-04-30"
$ws.Cells.Item(3, 11).Value = "蕭美琴"
$ws.Cells.Item(3, 12).Value = 981
$ws.Cells.Item(3, 13).Value = "tmpcd9a1"
$ws.Cells.Item(3, 14).Value = 114
